$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 16.67683263635534
$ws.Range("C2").Value = 10.38097268680016
$ws.Range("D2").Value = 5.193486321225011
$ws.Range("F2").Value = 30.80293401150741
$ws.Range("G2").Value = 3.63629011236276
$ws.Range("B3").Value = 16.03883890731553
$ws.Range("C3").Value = 9.771977088523654
$ws.Range("D3").Value = 5.229356604354764
$ws.Range("F3").Value = 30.47596236322454
$ws.Range("G3").Value = 3.640457525200651
$ws.Range("B4").Value = 15.63869190515078
$ws.Range("C4").Value = 9.379887392140407
$ws.Range("D4").Value = 5.252441166776948
$ws.Range("F4").Value = 30.28681697868037
$ws.Range("G4").Value = 3.643144890481822
$ws.Range("B5").Value = 15.473792629928
$ws.Range("C5").Value = 9.215678490133755
$ws.Range("D5").Value = 5.26211433653214
$ws.Range("F5").Value = 30.21272882776486
$ws.Range("G5").Value = 3.644272474614783
$ws.Range("B6").Value = 15.44630906621856
$ws.Range("C6").Value = 9.188148855616102
$ws.Range("D6").Value = 5.263736611182353
$ws.Range("F6").Value = 30.20060899069183
$ws.Range("G6").Value = 3.644461673599617
$ws.Range("B7").Value = 15.63647505955876
$ws.Range("C7").Value = 9.377690542187178
$ws.Range("D7").Value = 5.252570546079345
$ws.Range("F7").Value = 30.28580560813861
$ws.Range("G7").Value = 3.643159965863108
$ws.Range("B8").Value = 16.45875874611712
$ws.Range("C8").Value = 10.1748350284551
$ws.Range("D8").Value = 5.205634046286368
$ws.Range("F8").Value = 30.6878272153115
$ws.Range("G8").Value = 3.637700435376526
$ws.Range("B9").Value = 17.99356371313663
$ws.Range("C9").Value = 11.58919431802236
$ws.Range("D9").Value = 5.122027482849446
$ws.Range("F9").Value = 31.56492295315218
$ws.Range("G9").Value = 3.628008214777356
$ws.Range("B10").Value = 19.06109896774956
$ws.Range("C10").Value = 12.53280519796991
$ws.Range("D10").Value = 5.065779403189628
$ws.Range("F10").Value = 32.25831087965073
$ws.Range("G10").Value = 3.621496827577065
$ws.Range("B11").Value = 19.53151757435473
$ws.Range("C11").Value = 12.96794184944514
$ws.Range("D11").Value = 5.041323641164245
$ws.Range("F11").Value = 32.58318188672732
$ws.Range("G11").Value = 3.618665125162606
$ws.Range("B12").Value = 19.70731454004647
$ws.Range("C12").Value = 13.12992393975152
$ws.Range("D12").Value = 5.032226520169612
$ws.Range("F12").Value = 32.70745333457587
$ws.Range("G12").Value = 3.617611435901585
$ws.Range("B13").Value = 19.66955983072889
$ws.Range("C13").Value = 13.09520763776292
$ws.Range("D13").Value = 5.034178448625772
$ws.Range("F13").Value = 32.68063542001857
$ws.Range("G13").Value = 3.617837540908687
$ws.Range("B14").Value = 19.54602820631504
$ws.Range("C14").Value = 12.98134459510669
$ws.Range("D14").Value = 5.04057192786515
$ws.Range("F14").Value = 32.59338135465064
$ws.Range("G14").Value = 3.618578065189559
$ws.Range("B15").Value = 19.47005253577479
$ws.Range("C15").Value = 12.91110371715435
$ws.Range("D15").Value = 5.044509473867869
$ws.Range("F15").Value = 32.54009522601486
$ws.Range("G15").Value = 3.61903407803806
$ws.Range("B16").Value = 19.03003638395832
$ws.Range("C16").Value = 12.5057160241552
$ws.Range("D16").Value = 5.067400502505507
$ws.Range("F16").Value = 32.23726065842318
$ws.Range("G16").Value = 3.621684497998837
$ws.Range("B17").Value = 18.75608965042116
$ws.Range("C17").Value = 12.2659171622913
$ws.Range("D17").Value = 5.081733967701021
$ws.Range("F17").Value = 32.0538244486609
$ws.Range("G17").Value = 3.623343740175516
$ws.Range("B18").Value = 18.59709901212365
$ws.Range("C18").Value = 12.12598000202407
$ws.Range("D18").Value = 5.090084661355553
$ws.Range("F18").Value = 31.94921253749638
$ws.Range("G18").Value = 3.624310370323157
$ws.Range("B19").Value = 18.54302800440647
$ws.Range("C19").Value = 12.07825566910447
$ws.Range("D19").Value = 5.092930318737266
$ws.Range("F19").Value = 31.9139498427133
$ws.Range("G19").Value = 3.624639767409546
$ws.Range("B20").Value = 18.78540027465341
$ws.Range("C20").Value = 12.29165263865196
$ws.Range("D20").Value = 5.080197120296861
$ws.Range("F20").Value = 32.07325961663489
$ws.Range("G20").Value = 3.623165841241645
$ws.Range("B21").Value = 19.58237707602699
$ws.Range("C21").Value = 13.01489236077183
$ws.Range("D21").Value = 5.038689555703876
$ws.Range("F21").Value = 32.6189769446751
$ws.Range("G21").Value = 3.61836005109306
$ws.Range("B22").Value = 20.08954905876768
$ws.Range("C22").Value = 13.47929055191117
$ws.Range("D22").Value = 5.012516912996566
$ws.Range("F22").Value = 32.9828679852987
$ws.Range("G22").Value = 3.615327631248108
$ws.Range("B23").Value = 19.82016069319489
$ws.Range("C23").Value = 13.23346075310391
$ws.Range("D23").Value = 5.026398017574722
$ws.Range("F23").Value = 32.78802706565882
$ws.Range("G23").Value = 3.616936211491954
$ws.Range("B24").Value = 18.77215358621174
$ws.Range("C24").Value = 12.28002408541068
$ws.Range("D24").Value = 5.08089158585426
$ws.Range("F24").Value = 32.06447032764291
$ws.Range("G24").Value = 3.6232462297521
$ws.Range("B25").Value = 17.58808969580176
$ws.Range("C25").Value = 11.2230731593029
$ws.Range("D25").Value = 5.143738776083871
$ws.Range("F25").Value = 31.31864873390504
$ws.Range("G25").Value = 3.630522563279182
